$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.207.16"
$ws.Range("E2").Value = "  +1.55%  "

$ws.Range("D3").Value = "1.909.25"
$ws.Range("E3").Value = "  +1.83%  "

$ws.Range("E4").Value = "  -0.33%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "327.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.52%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.001"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.35%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4612"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.55%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3941"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.82%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "46.76"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.42%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07946"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +1.19%  "

$ws.Range("E11").Value = "  +1.06%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "22.34"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +2.54%  "

$ws.Range("D13").Value = "1.877.16"
$ws.Range("E13").Value = "  -0.87%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.103"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.52%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.769"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.98%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.06954"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.10%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "88.46"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.003"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.32%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.00001007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.23%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.13"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +1.76%  "

$ws.Range("E21").Value = "  -0.10%  "

$ws.Range("D22").Value = "29.234.56"
$ws.Range("E22").Value = "  +1.59%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.363"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.56%  "

$ws.Range("E24").Value = "  +0.51%  "

$ws.Range("D25").Value = "2.143.69"
$ws.Range("E25").Value = "  +0.67%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.061"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.76%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "156.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +2.38%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "19.47"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.14%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.129"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +5.75%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.997"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.25%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "118.84"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.15%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.09375"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.66%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.9266"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.80%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.351"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.88%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.354"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.98%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.266"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.84%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.207"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +4.67%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05839"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.29%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02107"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.82%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "7.956"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.23%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.39%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.5756"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.93%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.1804"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.71%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "9.965"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.42%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.267"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +4.67%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.04"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.84%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.5409"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.03%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.07077"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.94%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.877"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.70%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "2.550"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +5.84%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "112.95"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.72%  "
